# Auto-generated script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.866.90'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.894.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7880'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -3.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.35'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07210'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.27%  '
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7663'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.505'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.891.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.52'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.154'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.857.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007798'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.160.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.124'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +15.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1650'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.424'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.054'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.409'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.548'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("E32").Value = '  +4.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.115'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05560'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.273'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7429'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9974'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.615'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01921'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.777'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.143.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4425'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.863'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8500'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.84'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.52%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.881'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.980'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.474'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.027'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.15%  '
